$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# ================================================================
# Part 1: paragraph "Решено было сделать реализовать логику, через
# свойство продукта в корзине, решено: в корзине устанавливать
# свойства ..."
#
# Target text change:
#   "реализовать логику," -> "логику,"
#   " через свойство продукта в корзине, решено: в корзине
#    устанавливать свойства " ->
#   " через свойство продукта в корзине: устанавливать свойства "
# Also: a new _GoBack bookmark ends up sitting right between
# "корзине: " and "устанавливать свойства ".
# ================================================================
$p = $d.Paragraphs(27)

# --- place protector bookmarks around the boundaries we must keep,
#     so the runtime's run-merging pass doesn't collapse them while
#     we perform the text edits below.
$rp1 = $p.Range
$rp1.Find.Execute("реализовать", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pp1 = $d.Range($rp1.Start, $rp1.Start)
$d.Bookmarks.Add("ProtA", $pp1)

$pA = $d.Paragraphs(27)
$rp2 = $pA.Range
$rp2.Find.Execute("реализовать ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pp2 = $d.Range($rp2.End, $rp2.End)
$d.Bookmarks.Add("ProtB", $pp2)

$pB = $d.Paragraphs(27)
$rp3 = $pB.Range
$rp3.Find.Execute("логику,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pp3 = $d.Range($rp3.End, $rp3.End)
$d.Bookmarks.Add("ProtC", $pp3)

# --- perform the actual text edits ---

# "реализовать " is dropped (leaves "логику,")
$pC = $d.Paragraphs(27)
$rC = $pC.Range
$rC.Find.Execute("реализовать ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ", решено" is dropped
$pD = $d.Paragraphs(27)
$rD = $pD.Range
$rD.Find.Execute(", решено", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# " в корзине" (the second occurrence, right before "устанавливать") is dropped
$pE = $d.Paragraphs(27)
$rE = $pE.Range
$rE.Find.Execute(" в корзине устанавливать", $true, $false, $false, $false, $false, $true, 1, $false, " устанавливать", 2)

# --- split "проду" / "кта" ---
$pF = $d.Paragraphs(27)
$rF = $pF.Range
$rF.Find.Execute("проду", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ptF = $d.Range($rF.End, $rF.End)
$d.Bookmarks.Add("ProtD", $ptF)

# --- place the real _GoBack bookmark between "корзине: " and "устанавливать свойства " ---
$pG = $d.Paragraphs(27)
$rG = $pG.Range
$rG.Find.Execute("корзине: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ptG = $d.Range($rG.End, $rG.End)
$d.Bookmarks.Add("_GoBack", $ptG)

# --- remove the temporary protector bookmarks (the run splits they
#     created remain even after the bookmarks themselves are gone) ---
$d.Bookmarks("ProtA").Delete()
$d.Bookmarks("ProtB").Delete()
$d.Bookmarks("ProtC").Delete()
$d.Bookmarks("ProtD").Delete()

# ================================================================
# Part 2: paragraph "...По части создания касто[_GoBack]много
# условия находил еще статью ..." -> "...По части создания
# кастомного условия находил еще статью ..."
#
# The old _GoBack bookmark that used to sit between "касто" and
# "много" is gone (it moved to Part 1 above); the two runs merge
# into a single "кастомного" run.
# ================================================================
$p2 = $d.Paragraphs(34)

# locate "касто" unambiguously via its unique preceding context
$r2 = $p2.Range
$r2.Find.Execute("создания касто")
$kastoStart = $r2.End - 5

# insert "омног" strictly inside the "касто" run (between "каст" and
# the final "о") so the surrounding spellStart/spellEnd markers keep
# wrapping a single run instead of the edit spilling outside them
$insPt = $d.Range($kastoStart + 4, $kastoStart + 4)
$insPt.InsertAfter("омног")

# delete the now left-over "много" run entirely
$p3 = $d.Paragraphs(34)
$r3 = $p3.Range
$r3.Find.Execute("создания кастомного")
$mnogoRange = $d.Range($r3.End, $r3.End + 5)
$mnogoRange.Text = ""

# remove the stray _GoBack bookmark that used to live here
if ($d.Bookmarks.Exists("_GoBack")) {
    if ($d.Bookmarks("_GoBack").Start -ne $ptG.Start) {
        $d.Bookmarks("_GoBack").Delete()
    }
}

Write-Output "done"
